# Update gh-pages output generated at 456a3b4
# Applies numeric "want-to-go" count bumps across the four sheets and
# appends a new exhibition row to sheet "展览" (Worksheets.Item(1)).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws4 = $wb.Worksheets.Item(4)  # 全部类型

# ---------------------------------------------------------------
# Sheet 1 (展览): column F ("想去人数") updates
# ---------------------------------------------------------------
$ws1.Cells.Item(6, 6).Value = 830
$ws1.Cells.Item(7, 6).Value = 4236
$ws1.Cells.Item(9, 6).Value = 178
$ws1.Cells.Item(11, 6).Value = 6160
$ws1.Cells.Item(12, 6).Value = 6160
$ws1.Cells.Item(14, 6).Value = 468
$ws1.Cells.Item(15, 6).Value = 2356
$ws1.Cells.Item(19, 6).Value = 9289
$ws1.Cells.Item(21, 6).Value = 2507
$ws1.Cells.Item(23, 6).Value = 2328
$ws1.Cells.Item(24, 6).Value = 2473
$ws1.Cells.Item(29, 6).Value = 62
$ws1.Cells.Item(30, 6).Value = 336
$ws1.Cells.Item(32, 6).Value = 46
$ws1.Cells.Item(35, 6).Value = 76
$ws1.Cells.Item(36, 6).Value = 386
$ws1.Cells.Item(37, 6).Value = 1225
$ws1.Cells.Item(42, 6).Value = 1562
$ws1.Cells.Item(43, 6).Value = 2566
$ws1.Cells.Item(45, 6).Value = 933
$ws1.Cells.Item(48, 6).Value = 28

# ---------------------------------------------------------------
# Sheet 1 (展览): append new row 50
# ---------------------------------------------------------------
$ws1.Cells.Item(50, 1).Value = 49
$ws1.Cells.Item(50, 2).NumberFormat = "@"
$ws1.Cells.Item(50, 2).Value = "2024-11-16"
$ws1.Cells.Item(50, 3).Value = "北京·万游引力国潮动漫嘉年华S9"
$ws1.Cells.Item(50, 4).Value = "金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心"
$ws1.Cells.Item(50, 5).Value = "2024.11.16 10:00-11.17 17:00"
$ws1.Cells.Item(50, 6).Value = 0
$ws1.Cells.Item(50, 7).Value = 75
$ws1.Cells.Item(50, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90506"
$ws1.Cells.Item(50, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/Sp0BupUH1723136613917.jpeg"

# Fix up formatting introduced by the above edits so the new row matches
# the look & feel of the existing rows: column A keeps the bold/bordered
# "index" style, column B should not carry an explicit text format/style
# (it just needs to keep its value as literal text instead of becoming a
# date serial).
$ws1.Range("A49").Copy()
$ws1.Range("A50").PasteSpecial(-4122)

$ws1.Range("C49").Copy()
$ws1.Range("B50").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Sheet 2 (演出): column F updates
# ---------------------------------------------------------------
$ws2.Cells.Item(12, 6).Value = 151
$ws2.Cells.Item(22, 6).Value = 100

# ---------------------------------------------------------------
# Sheet 3 (本地生活): column F updates
# ---------------------------------------------------------------
$ws3.Cells.Item(3, 6).Value = 907

# ---------------------------------------------------------------
# Sheet 4 (全部类型): column F updates
# ---------------------------------------------------------------
$ws4.Cells.Item(4, 6).Value = 907
$ws4.Cells.Item(11, 6).Value = 830
$ws4.Cells.Item(12, 6).Value = 4236
$ws4.Cells.Item(14, 6).Value = 178
$ws4.Cells.Item(17, 6).Value = 6160
$ws4.Cells.Item(19, 6).Value = 2356
$ws4.Cells.Item(22, 6).Value = 9289
$ws4.Cells.Item(23, 6).Value = 151
$ws4.Cells.Item(25, 6).Value = 2507
$ws4.Cells.Item(27, 6).Value = 2473
$ws4.Cells.Item(32, 6).Value = 62
$ws4.Cells.Item(33, 6).Value = 336
$ws4.Cells.Item(34, 6).Value = 46
$ws4.Cells.Item(37, 6).Value = 76
$ws4.Cells.Item(38, 6).Value = 386
$ws4.Cells.Item(39, 6).Value = 1225
$ws4.Cells.Item(43, 6).Value = 2566
$ws4.Cells.Item(44, 6).Value = 933
$ws4.Cells.Item(48, 6).Value = 28
$ws4.Cells.Item(50, 6).Value = 100
$ws4.Cells.Item(51, 6).Value = 100
